$d = $word.ActiveDocument

function Rename-InlinePicture($range, $newName) {
    $ishapes = $range.InlineShapes
    if ($ishapes.Count -ge 1) {
        $shp = $ishapes.Item(1)
        $origAlt = $shp.AlternativeText
        $floating = $shp.ConvertToShape()
        $floating.Name = $newName
        $floating.AlternativeText = $origAlt
        [void]$floating.ConvertToInlineShape()
    }
}

$sec = $d.Sections.Item(1)

# Footer (primary) -> physical footer2.xml (docPr id="2"): image2.png -> image1.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image1.png"

# Footer (first page) -> physical footer1.xml (docPr id="3"): image2.png -> image1.png
Rename-InlinePicture $sec.Footers.Item(2).Range "image1.png"

# Header (first page) -> physical header1.xml (docPr id="1"): image1.jpg -> image2.jpg
Rename-InlinePicture $sec.Headers.Item(2).Range "image2.jpg"
